# Add a "field" forecast-argument column to the Pool sheet, and restore the
# author's active-tab/selection state (Pool tab active, Requirement no
# longer active).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pool")

# --- New "field" column (R) -------------------------------------------------
# Header
$ws.Range("R1").Value = "field"

# Data rows: write each boolean-looking value as a temporary helper formula
# (T("TRUE")/T("FALSE") forces a text result rather than the native Boolean
# type that a plain Value="TRUE" assignment would produce - matching the
# existing TRUE/FALSE text entries already used in columns K:M), then
# convert the whole block to static values with PasteSpecial so no formula
# is left behind.
$ws.Range("R2").Formula = "=T(""TRUE"")"
$ws.Range("R3").Formula = "=T(""FALSE"")"
$ws.Range("R4").Formula = "=T(""FALSE"")"
$ws.Range("R5").Formula = "=T(""FALSE"")"
$ws.Range("R6").Formula = "=T(""FALSE"")"

$ws.Range("R2:R6").Copy()
$ws.Range("R2:R6").PasteSpecial(-4163)
$excel.CutCopyMode = 0

# --- Active sheet / selection bookkeeping -----------------------------------
# Author moved focus from "Requirement" to "Pool" and left the cursor one
# row below the newly-added data (R7), scrolled so column C is visible.
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("R7").Select() | Out-Null
